$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '26.932.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '1.553.70'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +1.34%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.56%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '207.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +0.91%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  +0.85%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +0.51%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '21.77'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  +2.59%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  +1.47%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  +1.18%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'" + '0.0858'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  +0.55%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '1.775.08'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  +1.34%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '1.556.24'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  +1.76%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  +1.81%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +2.47%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '26.930.92'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  +0.68%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  +1.43%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '217.02'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  +2.15%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  +1.20%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '7.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  +0.26%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  +0.54%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +1.40%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '9.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  +1.51%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  +1.32%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '153.88'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  +1.72%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  +0.26%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '14.85'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  +0.63%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  +0.59%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  +1.39%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '0.0469'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +3.58%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +0.42%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +0.07%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '1.426.89'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +4.94%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  +3.63%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  +3.71%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  +0.74%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '2.30'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  +1.03%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  +0.18%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  +0.41%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '0.811'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  +1.45%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  +0.58%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '5.68'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -0.55%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'" + '0.990'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -0.33%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '2.28'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  +3.88%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '63.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  +2.14%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '1.74'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  +0.89%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'" + '1.689.95'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +1.47%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '86.19'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  +1.28%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'" + '0.0524'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  +4.46%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'" + '0.0₆0100'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  +3.53%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +1.65%  '
$ws.Range('E51').Style = 'Normal'
